# Trade #119 closed at 2026-02-17 09:28:48 - unknown UNKNOWN +0.000%
#
# Updates the "Summary", "Strategy Status", "All Trades" and "MarketMaking"
# sheets to reflect the newly closed trade #119.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Summary sheet — headline stats
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1200.69          # Current Capital
$summary.Range("B4").Value = 0.7              # Total P&L $
$summary.Range("B6").Value = 119              # Total Trades
$summary.Range("B7").Value = 54               # Winning Trades
$summary.Range("B9").Value = 45.38            # Win Rate %

# ---------------------------------------------------------------------
# Strategy Status sheet — MarketMaking row (row 4)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 100.69            # Capital
$status.Range("D4").Value = 119               # Trades
$status.Range("E4").Value = 0.7               # P&L $
$status.Range("F4").Value = 0.69              # P&L %
$status.Range("G4").Value = 45.38             # Win Rate %

# ---------------------------------------------------------------------
# Append the new trade row (#119) to both "All Trades" and
# "MarketMaking" sheets.
# ---------------------------------------------------------------------
$newRow = 120

function Add-TradeRow($ws, $row) {
    $ws.Cells.Item($row, 1).Value = 119
    # Leading apostrophe keeps this date-looking value as literal text
    # instead of being auto-converted to a date serial number.
    $ws.Cells.Item($row, 2).Value = "'2026-02-17"
    $ws.Cells.Item($row, 3).Value = "09:28:42"
    $ws.Cells.Item($row, 4).Value = "MarketMaking"
    $ws.Cells.Item($row, 5).Value = "DOWN"
    $ws.Cells.Item($row, 6).Value = 0.82
    $ws.Cells.Item($row, 7).Value = 0.826346
    $ws.Cells.Item($row, 8).Value = "CLOSED"
    $ws.Cells.Item($row, 9).Value = 0.7739
    $ws.Cells.Item($row, 10).Value = 0.01
    $ws.Cells.Item($row, 11).Value = 100.69
    $ws.Cells.Item($row, 12).Value = 0
    $ws.Cells.Item($row, 13).Value = 0
    $ws.Cells.Item($row, 14).Value = 0.6
    $ws.Cells.Item($row, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item($row, 16).Value = "early_exit"
    $ws.Cells.Item($row, 17).Value = 0.12
}

$allTrades = $wb.Worksheets.Item("All Trades")
Add-TradeRow $allTrades $newRow

$marketMaking = $wb.Worksheets.Item("MarketMaking")
Add-TradeRow $marketMaking $newRow
